# Apply updated "dSF" (column F) values on Sheet1.
# These reflect a data repull / recalculation for specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = 0
    5  = -1
    6  = 3
    7  = 3
    8  = -4
    14 = 0
    19 = 0
    26 = -2
    30 = 0
    34 = 1
    37 = 0
    42 = 0
    43 = 0
    45 = 2
    49 = 1
    50 = 2
    66 = -4
    68 = 0
    73 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
